$wb = $excel.ActiveWorkbook

# Sheet "上海" (Shanghai) - delete row 19 ("魅力足球")
$wsShanghai = $wb.Worksheets.Item("上海")
$wsShanghai.Rows(19).Delete()

# Sheet "卫视" (satellite TV) - delete column G ("GPTV")
$wsWeishi = $wb.Worksheets.Item("卫视")
$wsWeishi.Columns(7).Delete()
